$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("wickets")

$data = @(
  ,@(2, "9/51", "16/220")
  ,@(3, "8/71", "12/128")
  ,@(4, "7/42", "11/71")
  ,@(5, "10/74", "14/149")
  ,@(6, "8/24", "10/27")
  ,@(7, "8/15", "11/121")
  ,@(8, "7/37", "13/55")
  ,@(9, "7/51", "11/60")
  ,@(10, "9/83", "11/146")
  ,@(11, "9/127", "14/184")
  ,@(12, "9/52", "15/123")
  ,@(13, "7/87", "10/147")
  ,@(14, "8/84", "15/217")
  ,@(15, "7/119", "11/110")
  ,@(16, "7/59", "13/140")
  ,@(17, "8/45", "11/84")
  ,@(18, "8/50", "13/154")
  ,@(19, "7/37", "13/132")
  ,@(20, "8/34", "13/106")
  ,@(21, "7/22", "11/89")
  ,@(22, "7/76", "13/135")
  ,@(23, "8/58", "14/116")
  ,@(24, "7/87", "12/149")
  ,@(25, "7/83", "11/123")
  ,@(26, "7/71", "14/191")
  ,@(27, "8/71", "12/139")
  ,@(28, "8/43", "9/92")
  ,@(29, "7/64", "10/108")
  ,@(30, "8/61", "12/127")
  ,@(31, "7/87", "10/149")
  ,@(32, "5/30", "9/171")
  ,@(33, "6/23", "9/110")
  ,@(34, "8/38", "11/157")
  ,@(35, "8/31", "12/119")
  ,@(36, "7/74", "10/108")
  ,@(37, "8/51", "13/71")
  ,@(38, "6/30", "10/80")
  ,@(39, "6/54", "9/92")
  ,@(40, "8/97", "11/157")
  ,@(41, "7/98", "10/194")
  ,@(42, "7/77", "12/94")
  ,@(43, "6/56", "9/108")
  ,@(44, "7/37", "9/80")
  ,@(45, "6/50", "11/94")
  ,@(46, "6/65", "10/132")
  ,@(47, "7/39", "11/97")
  ,@(48, "8/92", "14/149")
  ,@(49, "7/72", "11/105")
  ,@(50, "7/61", "12/205")
  ,@(51, "8/71", "10/91")
  ,@(52, "8/79", "12/104")
  ,@(53, "7/44", "14/99")
  ,@(54, "8/86", "13/132")
  ,@(55, "9/56", "13/101")
  ,@(56, "8/41", "14/184")
  ,@(57, "6/73", "8/80")
  ,@(58, "7/46", "10/215")
  ,@(59, "6/26", "11/180")
  ,@(60, "7/39", "9/73")
  ,@(61, "6/42", "9/92")
  ,@(62, "7/38", "9/70")
  ,@(63, "7/12", "11/76")
  ,@(64, "5/58", "8/156")
  ,@(65, "6/21", "10/102")
  ,@(66, "6/48", "10/146")
  ,@(67, "7/48", "10/154")
  ,@(68, "6/54", "9/104")
  ,@(69, "7/27", "10/100")
  ,@(70, "7/40", "14/199")
  ,@(71, "6/73", "9/72")
  ,@(72, "7/36", "10/124")
  ,@(73, "7/112", "13/144")
  ,@(74, "6/67", "9/115")
  ,@(75, "8/87", "13/217")
  ,@(76, "8/108", "12/107")
  ,@(77, "8/164", "10/155")
  ,@(78, "7/54", "12/121")
  ,@(79, "7/40", "10/142")
  ,@(80, "6/46", "9/105")
)

foreach ($row in $data) {
  $r = $row[0]
  $hVal = $row[1]
  $iVal = $row[2]
  $hCell = $ws.Cells.Item($r, 8)
  $iCell = $ws.Cells.Item($r, 9)
  $hCell.NumberFormat = "@"
  $hCell.Value = $hVal
  $iCell.NumberFormat = "@"
  $iCell.Value = $iVal
}

$ws.Activate() | Out-Null
$ws.Range("R6").Select() | Out-Null